$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 47.5
$ws.Range("I6").Value = 56.666668
$ws.Range("K6").Value = 170.000004
$ws.Range("M6").Value = -58.00000399999999
$ws.Range("H19").Value = 1816
$ws.Range("I19").Value = 1793.6666
$ws.Range("J19").Value = 1849.5
$ws.Range("K19").Value = 1793.6666
$ws.Range("L19").Value = 1849.5
$ws.Range("M19").Value = -1618.6666
$ws.Range("N19").Value = -2199.5
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30584
$ws.Range("H51").Value = 10219.286
$ws.Range("J51").Value = 10264.708
$ws.Range("L51").Value = 10264.708
$ws.Range("N51").Value = -11232.708
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31068
$ws.Range("H92").Value = 691.25
$ws.Range("I92").Value = 404.78946
$ws.Range("K92").Value = 404.78946
$ws.Range("M92").Value = 843.21054
$ws.Range("H106").Value = 3115.9412
$ws.Range("I106").Value = 2467.5833
$ws.Range("J106").Value = 4672
$ws.Range("K106").Value = 2467.5833
$ws.Range("L106").Value = 4672
$ws.Range("M106").Value = -1836.5833
$ws.Range("N106").Value = -5934
$ws.Range("H129").Value = 2988.6428
$ws.Range("I129").Value = 1768.4
$ws.Range("K129").Value = 5305.200000000001
$ws.Range("M129").Value = -305.2000000000007

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 772215.7
$ws.Range("I45").Value = 911255
$ws.Range("K45").Value = 911255
$ws.Range("M45").Value = -910878
$ws.Range("H74").Value = 76931256
$ws.Range("I74").Value = 76931256
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 76931256
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -76930382
$ws.Range("H77").Value = 76931256
$ws.Range("I77").Value = 76931256
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 384656280
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -384651912
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2105.4614
$ws.Range("I20").Value = 1616.4445
$ws.Range("K20").Value = 1616.4445
$ws.Range("M20").Value = -1369.4445
$ws.Range("H86").Value = 3808.4285
$ws.Range("J86").Value = 3479.5
$ws.Range("L86").Value = 3479.5
$ws.Range("N86").Value = -5725.5
$ws.Range("H89").Value = 3808.4285
$ws.Range("J89").Value = 3479.5
$ws.Range("L89").Value = 17397.5
$ws.Range("N89").Value = -28629.5
$ws.Range("H105").Value = 1976.4445
$ws.Range("I105").Value = 1976.4445
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1976.4445
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -229.4445000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10282.6
$ws.Range("I22").Value = 16899.334
$ws.Range("J22").Value = 357.5
$ws.Range("K22").Value = 16899.334
$ws.Range("L22").Value = 357.5
$ws.Range("M22").Value = -16549.334
$ws.Range("N22").Value = -1057.5
$ws.Range("H31").Value = 14292.228
$ws.Range("I31").Value = 9719.888999999999
$ws.Range("J31").Value = 17457.691
$ws.Range("K31").Value = 9719.888999999999
$ws.Range("L31").Value = 17457.691
$ws.Range("M31").Value = -9424.888999999999
$ws.Range("N31").Value = -18047.691
$ws.Range("H34").Value = 14292.228
$ws.Range("I34").Value = 9719.888999999999
$ws.Range("J34").Value = 17457.691
$ws.Range("K34").Value = 9719.888999999999
$ws.Range("L34").Value = 17457.691
$ws.Range("M34").Value = -9517.888999999999
$ws.Range("N34").Value = -17861.691
$ws.Range("H141").Value = 324773.88
$ws.Range("J141").Value = 324773.88
$ws.Range("L141").Value = 324773.88
$ws.Range("N141").Value = -335133.88

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1389.875
$ws.Range("I129").Value = 1088.4286
$ws.Range("K129").Value = 3265.2858
$ws.Range("M129").Value = 1734.7142

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1458.9166
$ws.Range("I102").Value = 1278.6666
$ws.Range("J102").Value = 1999.6666
$ws.Range("K102").Value = 1278.6666
$ws.Range("L102").Value = 1999.6666
$ws.Range("M102").Value = 343.3334
$ws.Range("N102").Value = -5243.6666
$ws.Range("H113").Value = 204714.8
$ws.Range("I113").Value = 334924.66
$ws.Range("K113").Value = 334924.66
$ws.Range("M113").Value = -332754.66

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4433.875
$ws.Range("I16").Value = 4094.2
$ws.Range("K16").Value = 4094.2
$ws.Range("M16").Value = -3924.2
$ws.Range("H40").Value = 3216.276
$ws.Range("I40").Value = 3084.1482
$ws.Range("K40").Value = 3084.1482
$ws.Range("M40").Value = -2948.1482
$ws.Range("H46").Value = 1192.3334
$ws.Range("I46").Value = 1599.25
$ws.Range("J46").Value = 988.875
$ws.Range("K46").Value = 1599.25
$ws.Range("L46").Value = 988.875
$ws.Range("M46").Value = -1411.25
$ws.Range("N46").Value = -1364.875
$ws.Range("H74").Value = 79999
$ws.Range("J74").Value = 79999
$ws.Range("L74").Value = 79999
$ws.Range("N74").Value = -81995
$ws.Range("H77").Value = 79999
$ws.Range("J77").Value = 79999
$ws.Range("L77").Value = 239997
$ws.Range("N77").Value = -249981
$ws.Range("H93").Value = 4250
$ws.Range("I93").Value = 4250
$ws.Range("K93").Value = 4250
$ws.Range("M93").Value = -3002
$ws.Range("H136").Value = 1959.75
$ws.Range("I136").Value = 1731.3334
$ws.Range("J136").Value = 1992.381
$ws.Range("K136").Value = 5194.0002
$ws.Range("L136").Value = 5977.143
$ws.Range("M136").Value = -2644.0002
$ws.Range("N136").Value = -11077.143

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8469.1
$ws.Range("I4").Value = 15239.6
$ws.Range("J4").Value = 1698.6
$ws.Range("K4").Value = 15239.6
$ws.Range("L4").Value = 1698.6
$ws.Range("M4").Value = -15126.6
$ws.Range("N4").Value = -1924.6
$ws.Range("H6").Value = 25000
$ws.Range("I6").Value = 25000
$ws.Range("K6").Value = 25000
$ws.Range("M6").Value = -24885
$ws.Range("H74").Value = 16165.5
$ws.Range("J74").Value = 15999
$ws.Range("L74").Value = 15999
$ws.Range("N74").Value = -17871
$ws.Range("H77").Value = 16165.5
$ws.Range("J77").Value = 15999
$ws.Range("L77").Value = 47997
$ws.Range("N77").Value = -57357
$ws.Range("H82").Value = 41666.5
$ws.Range("J82").Value = 41666.5
$ws.Range("L82").Value = 41666.5
$ws.Range("N82").Value = -42432.5
$ws.Range("H85").Value = 41666.5
$ws.Range("J85").Value = 41666.5
$ws.Range("L85").Value = 41666.5
$ws.Range("N85").Value = -44318.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

